$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 40, shifting existing rows 40:46 down to 41:47
$ws.Rows.Item(40).Insert()

# Copy the style of the Fecha (date) cell from the row below (now row 41) to the new row 40
$ws.Cells.Item(41, 4).Copy()
$ws.Cells.Item(40, 4).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new weekly record in row 40
$ws.Cells.Item(40, 1).Value = 10
$ws.Cells.Item(40, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(40, 3).Value = "La Araucanía"
$ws.Cells.Item(40, 4).Value = 44449
$ws.Cells.Item(40, 5).Value = 9
$ws.Cells.Item(40, 6).Value = 100112035
$ws.Cells.Item(40, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(40, 8).Value = "Sin especificar"
$ws.Cells.Item(40, 9).Value = "Primera"
$ws.Cells.Item(40, 10).Value = 12
$ws.Cells.Item(40, 11).Value = 25000
$ws.Cells.Item(40, 12).Value = 25000
$ws.Cells.Item(40, 13).Value = 25000
$ws.Cells.Item(40, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(40, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(40, 16).Value = 2500
$ws.Cells.Item(40, 17).Value = 10
$ws.Cells.Item(40, 18).Value = "Hortaliza"
